$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "35.805.63"
$ws.Range("E2").Value = "  -2.56%  "

# Row 3
$ws.Range("D3").Value = "1.990.01"
$ws.Range("E3").Value = "  -3.46%  "

# Row 4
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.41%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.640"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.31%  "

# Row 7
$ws.Range("E7").Value = "  +7.15%  "

# Row 8
$ws.Range("E8").Value = "  +0.06%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "59.01"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.65%  "

# Row 10
$ws.Range("E10").Value = "  -0.81%  "

# Row 11
$ws.Range("E11").Value = "  -1.39%  "

# Row 12
$ws.Range("E12").Value = "  -2.38%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.958"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.31%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.59%  "

# Row 15
$ws.Range("D15").Value = "2.278.28"
$ws.Range("E15").Value = "  -3.49%  "

# Row 16
$ws.Range("E16").Value = "  -3.21%  "

# Row 17
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.84"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +8.09%  "

# Row 18
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "1.979.79"
$ws.Range("E18").Value = "  -3.91%  "

# Row 19
$ws.Range("D19").Value = "35.776.07"
$ws.Range("E19").Value = "  -2.33%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.83%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0853"
$ws.Range("E21").Value = "  -1.66%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.48%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.70"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.20%  "

# Row 24
$ws.Range("E24").Value = "  +0.07%  "

# Row 25
$ws.Range("E25").Value = "  +15.05%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.65%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.79%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "164.93"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.72%  "

# Row 29
$ws.Range("E29").Value = "  -4.02%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.119"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.45%  "

# Row 31
$ws.Range("E31").Value = "  -3.08%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.61%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0984"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +15.71%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0604"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.68%  "

# Row 35
$ws.Range("E35").Value = "  +10.88%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.41"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.85%  "

# Row 37
$ws.Range("E37").Value = "  +0.07%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.80"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.40%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.76"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +12.80%  "

# Row 40
$ws.Range("E40").Value = "  -1.52%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0961"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.57%  "

# Row 42
$ws.Range("B42").Value = "HuobiToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.87"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.98%  "

# Row 43
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0214"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.15%  "

# Row 44
$ws.Range("E44").Value = "  -1.01%  "

# Row 45
$ws.Range("E45").Value = "  -1.20%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.47"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.28%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.79"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.98%  "

# Row 48
$ws.Range("D48").Value = "1.374.53"
$ws.Range("E48").Value = "  -3.50%  "

# Row 49
$ws.Range("E49").Value = "  -1.30%  "

# Row 50
$ws.Range("E50").Value = "  +1.41%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "47.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.74%  "
